$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 94, shifting existing rows 94:227 down to 95:228.
$ws.Rows(94).Insert()

# Populate the new row 94 with the new weekly price entry (same template
# values as the former row 94, but a new date and a corrected
# commercialization unit).
$ws.Range("A94").Value = 5
$ws.Range("B94").Value = "Macroferia Regional de Talca"
$ws.Range("C94").Value = "Maule"
$ws.Range("D94").Value = 44571
$ws.Range("E94").Value = 7
$ws.Range("F94").Value = 100112003
$ws.Range("G94").Value = "Ajo"
$ws.Range("H94").Value = "Chino"
$ws.Range("I94").Value = "Primera"
$ws.Range("J94").Value = 200
$ws.Range("K94").Value = 20000
$ws.Range("L94").Value = 20000
$ws.Range("M94").Value = 20000
$ws.Range("N94").Value = "`$/caja 10 kilos"
$ws.Range("O94").Value = "China"
$ws.Range("P94").Value = 2000
$ws.Range("Q94").Value = 10
$ws.Range("R94").Value = "Hortaliza"
